$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update shortname value (shared by both sheets) to include hyphen
$ws1.Range("B1").Value = "295-MS-EPP-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"
$ws2.Range("B1").Value = "295-MS-EPP-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

# Switch selection on input sheet to B1
$ws1.Range("B1").Select()

# Activate output sheet and select B1
$ws2.Activate()
$ws2.Range("B1").Select()
